$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.752.85"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "2.126.18"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5278"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4578"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09124"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "2.137.84"
$ws.Range("E13").Value = "  +2.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.876"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.143"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.015"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06724"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("E20").Value = "  +2.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.013"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.391"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").Value = "30.815.08"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.367"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").Value = "2.368.47"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.26%  "
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.57"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.210"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1080"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.659"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.400"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.956"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.64"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.927"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02682"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06895"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2333"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.99%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6938"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.265"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.26%  "
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.320"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("E47").Value = "  +19.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.709"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.261"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07314"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.74%  "
